# Atualização de bases das ligas, do dia: 02-05-2024 às 20:28
# Swap the data (columns B through AB) between row 25 <-> row 26,
# and between row 85 <-> row 86. Column A (the running index) is left
# untouched on each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData {
    param($row1, $row2)

    # Columns B (2) through AB (28)
    $firstCol = 2
    $lastCol = 28

    $range1 = $ws.Range($ws.Cells.Item($row1, $firstCol), $ws.Cells.Item($row1, $lastCol))
    $range2 = $ws.Range($ws.Cells.Item($row2, $firstCol), $ws.Cells.Item($row2, $lastCol))

    $values1 = $range1.Value2
    $values2 = $range2.Value2

    $range1.Value2 = $values2
    $range2.Value2 = $values1
}

Swap-RowData 25 26
Swap-RowData 85 86
